$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look like plain numbers need an explicit Text
# number format first, otherwise Excel would silently convert them to
# numeric values (dropping things like trailing zeros).
$textCells = @(
    "D5", "D6", "D13", "D14", "D15", "D19", "D21", "D25", "D26", "D27", "D30", "D31", "D32", "D36", "D38", "D39", "D41", "D42", "D46", "D47", "D48", "D49"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "63.282.79"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "2.485.27"
$ws.Range("E3").Value = "  +3.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").Value = "577.97"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6
$ws.Range("D6").Value = "147.11"
$ws.Range("E6").Value = "  +0.95%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("D9").Value = "2.481.99"
$ws.Range("E9").Value = "  +1.89%  "

# Row 10
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("E11").Value = "  +1.77%  "

# Row 12
$ws.Range("E12").Value = "  +0.52%  "

# Row 13
$ws.Range("D13").Value = "0.354"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").Value = "28.66"
$ws.Range("E14").Value = "  +4.15%  "

# Row 15
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16
$ws.Range("D16").Value = "2.933.79"

# Row 17
$ws.Range("D17").Value = "63.202.78"
$ws.Range("E17").Value = "  +1.42%  "

# Row 18
$ws.Range("D18").Value = "2.480.56"
$ws.Range("E18").Value = "  +1.72%  "

# Row 19
$ws.Range("D19").Value = "8.23"
$ws.Range("E19").Value = "  +3.98%  "

# Row 20
$ws.Range("E20").Value = "  +0.99%  "

# Row 21
$ws.Range("D21").Value = "330.08"
$ws.Range("E21").Value = "  +0.88%  "

# Row 22
$ws.Range("E22").Value = "  +10.51%  "

# Row 24
$ws.Range("E24").Value = "  +0.18%  "

# Row 25
$ws.Range("D25").Value = "66.33"
$ws.Range("E25").Value = "  +1.11%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "9.85"
$ws.Range("E26").Value = "  +16.44%  "

# Row 27
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "669.80"
$ws.Range("E27").Value = "  +6.25%  "

# Row 28
$ws.Range("E28").Value = "  +1.77%  "

# Row 29
$ws.Range("D29").Value = "2.602.91"

# Row 30
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -9.61%  "

# Row 31
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +4.44%  "

# Row 32
$ws.Range("D32").Value = "8.10"
$ws.Range("E32").Value = "  -1.02%  "

# Row 33
$ws.Range("E33").Value = "  -0.40%  "

# Row 34
$ws.Range("E34").Value = "  -3.25%  "

# Row 35
$ws.Range("E35").Value = "  +4.39%  "

# Row 36
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  +0.29%  "

# Row 37
$ws.Range("E37").Value = "  +0.96%  "

# Row 38
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").Value = "  +2.05%  "

# Row 39
$ws.Range("D39").Value = "0.372"
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("E40").Value = "  +0.87%  "

# Row 41
$ws.Range("D41").Value = "150.70"
$ws.Range("E41").Value = "  -0.84%  "

# Row 42
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("E43").Value = "  +0.65%  "

# Row 45
$ws.Range("D45").Value = "0.0₆0314"
$ws.Range("E45").Value = "  -10.85%  "

# Row 46
$ws.Range("D46").Value = "156.25"
$ws.Range("E46").Value = "  +7.87%  "

# Row 47
$ws.Range("D47").Value = "15.25"
$ws.Range("E47").Value = "  +3.41%  "

# Row 48
$ws.Range("D48").Value = "3.62"
$ws.Range("E48").Value = "  +0.61%  "

# Row 49
$ws.Range("D49").Value = "20.55"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("E50").Value = "  +1.58%  "

# Row 51
$ws.Range("E51").Value = "  -0.01%  "
